$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 637
$ws.Range("F3").Value = 701
$ws.Range("F4").Value = 944
$ws.Range("F5").Value = 719
$ws.Range("F8").Value = 599
$ws.Range("F9").Value = 130
$ws.Range("F12").Value = 383
$ws.Range("F13").Value = 506
$ws.Range("F15").Value = 12
$ws.Range("F16").Value = 472
$ws.Range("F17").Value = 353
$ws.Range("F18").Value = 54
$ws.Range("F20").Value = 554
$ws.Range("F22").Value = 575
$ws.Range("F24").Value = 755
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 20
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 637
$ws.Range("F7").Value = 701
$ws.Range("F8").Value = 944
$ws.Range("F9").Value = 719
$ws.Range("F12").Value = 599
$ws.Range("F13").Value = 130
$ws.Range("F17").Value = 20
$ws.Range("F18").Value = 383
$ws.Range("F19").Value = 506
$ws.Range("F22").Value = 12
$ws.Range("F23").Value = 472
$ws.Range("F25").Value = 353
$ws.Range("F26").Value = 54
$ws.Range("F30").Value = 554
$ws.Range("F36").Value = 575
$ws.Range("F38").Value = 755
